$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the header labels
$ws.Range("B1").Value = "Primary Stat"
$ws.Range("C1").Value = "Image Name"

# Fix the "Abbadon" -> "Abaddon" spelling correction
$ws.Range("A2").Value = "Abaddon"
$ws.Range("C2").Value = "Abaddon"

# Update the view: scroll back to top and select C3
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C3").Select()
